$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6312
$ws.Range("I11").Value = 6312
$ws.Range("K11").Value = 6312
$ws.Range("M11").Value = -6172

$ws.Range("H33").Value = 296.48648
$ws.Range("I33").Value = 315.36667
$ws.Range("J33").Value = 215.57143
$ws.Range("K33").Value = 315.36667
$ws.Range("L33").Value = 215.57143
$ws.Range("M33").Value = -86.36667
$ws.Range("N33").Value = -673.57143

$ws.Range("H40").Value = 1123
$ws.Range("I40").Value = 1038.1818
$ws.Range("J40").Value = 1434
$ws.Range("K40").Value = 1038.1818
$ws.Range("L40").Value = 1434
$ws.Range("M40").Value = -863.1818000000001
$ws.Range("N40").Value = -1784

$ws.Range("H64").Value = 42183.23
$ws.Range("I64").Value = 4042.4443
$ws.Range("J64").Value = 128000
$ws.Range("K64").Value = 4042.4443
$ws.Range("L64").Value = 128000
$ws.Range("M64").Value = -3794.4443
$ws.Range("N64").Value = -128496

$ws.Range("H67").Value = 42183.23
$ws.Range("I67").Value = 4042.4443
$ws.Range("J67").Value = 128000
$ws.Range("K67").Value = 4042.4443
$ws.Range("L67").Value = 128000
$ws.Range("M67").Value = -3184.4443
$ws.Range("N67").Value = -129716

$ws.Range("H121").Value = 271.53845
$ws.Range("J121").Value = 262.4
$ws.Range("L121").Value = 787.1999999999999
$ws.Range("N121").Value = -4281.2

$ws.Range("H135").Value = 1324.5883
$ws.Range("I135").Value = 625.9259
$ws.Range("J135").Value = 4019.4285
$ws.Range("K135").Value = 5633.3331
$ws.Range("L135").Value = 36174.8565
$ws.Range("M135").Value = -3098.3331
$ws.Range("N135").Value = -41244.8565

$ws.Range("H137").Value = 3248.8909
$ws.Range("I137").Value = 3842.8286
$ws.Range("J137").Value = 2209.5
$ws.Range("K137").Value = 11528.4858
$ws.Range("L137").Value = 6628.5
$ws.Range("M137").Value = -8978.485799999999
$ws.Range("N137").Value = -11728.5

$ws.Range("H138").Value = 1257740.2
$ws.Range("I138").Value = 5386.273
$ws.Range("J138").Value = 1557216.1
$ws.Range("K138").Value = 16158.819
$ws.Range("L138").Value = 4671648.300000001
$ws.Range("M138").Value = -11018.819
$ws.Range("N138").Value = -4681928.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 65175.812
$ws.Range("J45").Value = 2711.4
$ws.Range("L45").Value = 2711.4
$ws.Range("N45").Value = -3465.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 202.625
$ws.Range("I22").Value = 203.14285
$ws.Range("K22").Value = 203.14285
$ws.Range("M22").Value = -30.14285000000001

$ws.Range("H103").Value = 23266.8
$ws.Range("J103").Value = 23266.8
$ws.Range("L103").Value = 23266.8
$ws.Range("N103").Value = -25610.8

$ws.Range("H134").Value = 35836.625
$ws.Range("I134").Value = 43295.08
$ws.Range("J134").Value = 3516.6667
$ws.Range("K134").Value = 129885.24
$ws.Range("L134").Value = 10550.0001
$ws.Range("M134").Value = -127350.24
$ws.Range("N134").Value = -15620.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3031851.2
$ws.Range("I31").Value = 1280.0312
$ws.Range("K31").Value = 1280.0312
$ws.Range("M31").Value = -985.0311999999999

$ws.Range("H34").Value = 3031851.2
$ws.Range("I34").Value = 1280.0312
$ws.Range("K34").Value = 1280.0312
$ws.Range("M34").Value = -1078.0312

$ws.Range("H68").Value = 20557.143
$ws.Range("I68").Value = 20000
$ws.Range("J68").Value = 20650
$ws.Range("K68").Value = 20000
$ws.Range("L68").Value = 20650
$ws.Range("M68").Value = -19251
$ws.Range("N68").Value = -22148

$ws.Range("H71").Value = 20557.143
$ws.Range("I71").Value = 20000
$ws.Range("J71").Value = 20650
$ws.Range("K71").Value = 60000
$ws.Range("L71").Value = 61950
$ws.Range("M71").Value = -56256
$ws.Range("N71").Value = -69438

$ws.Range("H81").Value = 50327.25
$ws.Range("J81").Value = 50327.25
$ws.Range("L81").Value = 50327.25
$ws.Range("N81").Value = -52323.25

$ws.Range("H84").Value = 50327.25
$ws.Range("J84").Value = 50327.25
$ws.Range("L84").Value = 150981.75
$ws.Range("N84").Value = -160965.75

$ws.Range("H105").Value = 989.7
$ws.Range("I105").Value = 899.875
$ws.Range("J105").Value = 1349
$ws.Range("K105").Value = 899.875
$ws.Range("L105").Value = 1349
$ws.Range("M105").Value = 847.125
$ws.Range("N105").Value = -4843

$ws.Range("H134").Value = 2133.6177
$ws.Range("I134").Value = 2060.7407
$ws.Range("J134").Value = 2414.7144
$ws.Range("K134").Value = 6182.222099999999
$ws.Range("L134").Value = 7244.1432
$ws.Range("M134").Value = -3647.222099999999
$ws.Range("N134").Value = -12314.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1239.5862
$ws.Range("I68").Value = 860.5454999999999
$ws.Range("J68").Value = 2430.8572
$ws.Range("K68").Value = 2581.6365
$ws.Range("L68").Value = 7292.571599999999
$ws.Range("M68").Value = -1770.6365
$ws.Range("N68").Value = -8914.571599999999

$ws.Range("H71").Value = 1239.5862
$ws.Range("I71").Value = 860.5454999999999
$ws.Range("J71").Value = 2430.8572
$ws.Range("K71").Value = 7744.9095
$ws.Range("L71").Value = 21877.7148
$ws.Range("M71").Value = -3688.9095
$ws.Range("N71").Value = -29989.7148

$ws.Range("H107").Value = 21710.105
$ws.Range("I107").Value = 19387.32
$ws.Range("J107").Value = 24641.238
$ws.Range("K107").Value = 58161.96
$ws.Range("L107").Value = 73923.71400000001
$ws.Range("M107").Value = -56241.96
$ws.Range("N107").Value = -77763.71400000001

$ws.Range("H137").Value = 49920.637
$ws.Range("I137").Value = 2513.5
$ws.Range("J137").Value = 132883.12
$ws.Range("K137").Value = 7540.5
$ws.Range("L137").Value = 398649.36
$ws.Range("M137").Value = -2440.5
$ws.Range("N137").Value = -408849.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 122777.4
$ws.Range("I80").Value = 4753.6
$ws.Range("J80").Value = 240801.2
$ws.Range("K80").Value = 4753.6
$ws.Range("L80").Value = 240801.2
$ws.Range("M80").Value = -3755.6
$ws.Range("N80").Value = -242797.2

$ws.Range("H83").Value = 122777.4
$ws.Range("I83").Value = 4753.6
$ws.Range("J83").Value = 240801.2
$ws.Range("K83").Value = 23768
$ws.Range("L83").Value = 1204006
$ws.Range("M83").Value = -18776
$ws.Range("N83").Value = -1213990

$ws.Range("H132").Value = 17243320
$ws.Range("I132").Value = 50001980
$ws.Range("J132").Value = 1918.1842
$ws.Range("K132").Value = 150005940
$ws.Range("L132").Value = 5754.5526
$ws.Range("M132").Value = -150003410
$ws.Range("N132").Value = -10814.5526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 368.66666
$ws.Range("I22").Value = 312.25
$ws.Range("J22").Value = 820
$ws.Range("K22").Value = 312.25
$ws.Range("L22").Value = 820
$ws.Range("M22").Value = -17.25
$ws.Range("N22").Value = -1410

$ws.Range("H27").Value = 368.66666
$ws.Range("I27").Value = 312.25
$ws.Range("J27").Value = 820
$ws.Range("K27").Value = 312.25
$ws.Range("L27").Value = 820
$ws.Range("M27").Value = -205.25
$ws.Range("N27").Value = -1034

$ws.Range("H122").Value = 19767
$ws.Range("I122").Value = 51752
$ws.Range("J122").Value = 3774.5
$ws.Range("K122").Value = 155256
$ws.Range("L122").Value = 11323.5
$ws.Range("M122").Value = -152806
$ws.Range("N122").Value = -16223.5

$ws.Range("H132").Value = 5006.698
$ws.Range("I132").Value = 5478.027
$ws.Range("J132").Value = 3916.75
$ws.Range("K132").Value = 16434.081
$ws.Range("L132").Value = 11750.25
$ws.Range("M132").Value = -13904.081
$ws.Range("N132").Value = -16810.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3657
$ws.Range("I81").Value = 1100
$ws.Range("J81").Value = 4679.8
$ws.Range("K81").Value = 2200
$ws.Range("L81").Value = 9359.6
$ws.Range("M81").Value = -1139
$ws.Range("N81").Value = -11481.6

$ws.Range("H84").Value = 3657
$ws.Range("I84").Value = 1100
$ws.Range("J84").Value = 4679.8
$ws.Range("K84").Value = 11000
$ws.Range("L84").Value = 46798
$ws.Range("M84").Value = -5696
$ws.Range("N84").Value = -57406

$ws.Range("H132").Value = 7007011
$ws.Range("I132").Value = 8719436
$ws.Range("J132").Value = 1635.6364
$ws.Range("K132").Value = 26158308
$ws.Range("L132").Value = 4906.9092
$ws.Range("M132").Value = -26155778
$ws.Range("N132").Value = -9966.9092

$ws.Range("H136").Value = 8475656
$ws.Range("I136").Value = 19637.455
$ws.Range("J136").Value = 23978356
$ws.Range("K136").Value = 58912.36500000001
$ws.Range("L136").Value = 71935068
$ws.Range("M136").Value = -56362.36500000001
$ws.Range("N136").Value = -71940168
